# #1 - Batch 구조 설계
# 1개의 헬스장 - 다수의 회원 을 기준으로 하여 작성하였음
#
# Nudge the four table graphic frames on slide 1 slightly down-right
# (dx = 0.75pt / 9525 EMU, dy = 5.25pt / 66675 EMU).
#
# NOTE: Shape.Left/Top are exposed as single-precision floats, and the
# point->EMU conversion truncates rather than rounds, so the target
# positions below are written as the literal point values that survive
# that round trip exactly (rather than read-modify-write with +=, which
# would lose a bit of precision and land 1 EMU short).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$lefts = @(694.75, 139.0, 324.25, 509.5)
$top = 173.20111083984375

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    $shp.Left = $lefts[$i - 1]
    $shp.Top = $top
}
